$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) keeps its existing bold/centered header style ---
$headers = @("code", "product_id", "product_name", "price", "promo_price", "date_from", "date_to")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Data rows (rows 2-9) ---
$data = @(
    @("745883670185", "6080", "Uchwyt samochodowy Belkin Car Cup Mount (do 83mm), czarny", "169.00", "149.00", "25-05-2023", "31-12-2027"),
    @("5901571174976", "8272", "Folia ochronna 3mk Arc Microsoft Band 2, 3 szt.", "29.00", "11.60", "10-04-2024", "31-12-2027"),
    @("8809565309809", "9847", "Etui wodoszczelne Spigen Velo A600 IPX8  do telefonu max 16,3 x 9cm, białe", "69.00", "41.40", "10-04-2024", "31-12-2027"),
    @("745883734108", "9876", "Adapter Belkin MFi 3.5mm Audio + Charge RockStar, białe", "205.00", "159.00", "25-05-2023", "31-12-2027"),
    @("5904538092801", "10855", "Szkło hartowane na aparat Mocolo TG+ iPhone Xs/X", "19.00", "15.20", "12-02-2024", "31-12-2027"),
    @("5904538092825", "10858", "Szkło hartowane na aparat Mocolo TG+ iPhone Xs Max", "19.00", "7.60", "08-01-2024", "31-12-2027"),
    @("6953156257344", "11566", "Uchwyt samochodowy na telefon Baseus Gravity Type SUYL-01, czarny", "49.00", "35.00", "26-01-2025", "01-12-2027"),
    @("8809613769999", "11619", "Etui Spigen Liquid Air Galaxy A9 2018, matowe czarne", "49.00", "1.99", "24-10-2022", "31-12-2027")
)

# Force text storage so numeric-looking codes/prices and dd-mm-yyyy dates
# are kept verbatim as strings instead of being auto-converted by Excel
# into numbers / serial dates.
$dataRange = $ws.Range("A2:G9")
$dataRange.NumberFormat = "@"

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Restore the default (unstyled) look for the data cells now that the
# values are safely stored as text.
$dataRange.Style = "Normal"

# --- Remove the now-unused condition_type / product_id(dup) / stocks columns ---
$ws.Columns("H:J").Delete()
